# Commit: "Delete Create Address file"
#
# The workbook gains a new worksheet "NewShppingAddress" inserted right
# after "BillingAddress". It is (mostly) a copy of the BillingAddress row,
# but with a different contact (Jihad Mohammad / Integrant) and an extra
# phone-number column (J). BillingAddress itself keeps its data but its
# view/selection and column widths are tweaked, and it stops being the
# active/selected tab (the new sheet becomes active instead).

$wb = $excel.ActiveWorkbook

$billing = $wb.Worksheets.Item("BillingAddress")

# ---------------------------------------------------------------------
# 1. Tweak BillingAddress first: new column width + new selection
#    (done before activating the new sheet, since the last-selected
#    sheet becomes the workbook's active tab)
# ---------------------------------------------------------------------
$billing.Columns.Item(10).ColumnWidth = 15.999999999999998
$billing.Range("A2:I2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new worksheet right after BillingAddress
# ---------------------------------------------------------------------
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $billing)
$newWs.Name = "NewShppingAddress"

# ---------------------------------------------------------------------
# 3. Fill in the new sheet's data (row 2, columns A:J)
# ---------------------------------------------------------------------
$newWs.Range("A2").Value = "Integrant"
$newWs.Range("B2").Value = "Jihad"
$newWs.Range("C2").Value = "Mohammad"
$newWs.Range("D2").Value = "JihadMohamed@Gmail.com"
$newWs.Range("E2").Value = 12
$newWs.Range("F2").Value = 566
$newWs.Range("G2").Value = 11411
$newWs.Range("H2").Value = "FR"
$newWs.Range("I2").Value = "France"
$newWs.Range("J2").Value = 1257897445

# Hyperlink + hyperlink style (matches the style used elsewhere for the
# same kind of "email" cell)
$newWs.Hyperlinks.Add($newWs.Range("D2"), "mailto:JihadMohamed@Gmail.com") | Out-Null
$newWs.Range("D2").Style = "Hyperlink"

# Column widths for the new sheet
$newWs.Columns.Item(3).ColumnWidth = 12.5
$newWs.Columns.Item(4).ColumnWidth = 16.666666666666668
$newWs.Columns.Item(10).ColumnWidth = 17.666666666666668

# Selection / active view for the new (now active) sheet -- this must be
# the last selection made so the new sheet ends up as the active tab.
$newWs.Range("J5").Select() | Out-Null
